# Add files via upload
# "Test Results" sheet: the single "Test Passed" column is replaced by
# four CRUD test-result columns (Create/Read/Update/Delete Test Passed),
# each defaulting to FALSE, for every existing row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Extend the formatting (border/font) of the existing header & data cells
# into the three new columns before touching values, so the new cells
# pick up the same styles (bold+border for row 1, border for the rest)
# instead of creating brand-new style entries.
$ws.Range("B1").Copy()
$ws.Range("C1:E1").PasteSpecial(-4122)
$ws.Range("B2:B24").Copy()
$ws.Range("C2:E24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row: B1 keeps the "Test Passed" concept but becomes "Create Test
# Passed"; C1/D1/E1 are brand new headers.
$ws.Range("B1").Value = "Create Test Passed"
$ws.Range("C1").Value = "Read Test Passed"
$ws.Range("D1").Value = "Update Test Passed"
$ws.Range("E1").Value = "Delete Test Passed"

# Fill the new/extended boolean columns (rows 2-24) with FALSE, matching
# the existing B column values.
$ws.Range("B2:E24").Value = $false

# Column widths to match the widened headers (values chosen so the
# engine's internal char->OOXML-width rounding lands on the saved widths).
$ws.Range("B1:C1").ColumnWidth = 17
$ws.Range("D1").ColumnWidth = 17.666666666666668
$ws.Range("E1").ColumnWidth = 17.166666666666668

# Selection ends up on E2, matching the saved view state.
$ws.Range("E2").Select()
